# Update the "Sexo" column (G) values:
#   Masculino -> Hombre
#   Femenino  -> Mujer
# Also restore the active-cell selection to match the target workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2
    if ($val -eq "Masculino") {
        $cell.Value2 = "Hombre"
    } elseif ($val -eq "Femenino") {
        $cell.Value2 = "Mujer"
    }
}

# Restore selection to V15 as in the target file.
$ws.Range("V15").Select()
